$wb = $excel.ActiveWorkbook

# --- Sheet "book-list": add a new row for book #4 ---
$booklist = $wb.Worksheets.Item("book-list")
$booklist.Range("A5").Value = 4
$booklist.Range("B5").Value = "The Sensory Order"
$booklist.Range("C5").Value = 1952
$booklist.Range("D5").Value = "Friedrich A. Hayek"

# --- Sheet "reading-data": mark rows 71-95 (book column G) as book #4 ---
$readingdata = $wb.Worksheets.Item("reading-data")
$readingdata.Range("G71:G95").Value = 4

# --- Update the active selection/view on reading-data to match the new state ---
$readingdata.Activate()
$readingdata.Range("G72:G95").Select()
